# LoginData.xlsx edit: update row 2 credentials on the "Login" sheet,
# remove the now-stale hyperlink on C2, and move the active selection to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Remove the hyperlink that currently sits on C2 (mailto:Tetherfi@930).
# $ws.Hyperlinks.Delete(range) / range.Hyperlinks.Delete() both clear every
# hyperlink on the sheet here, so find and delete the single Hyperlink
# object whose Range is C2 instead.
$target = $null
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$2') {
        $target = $hl
    }
}
if ($target -ne $null) {
    $target.Delete()
}

# New credentials for row 2. C2 is written before B2 so the new shared
# strings land in the same table order as the target file (P@ssw0rd@123
# at index 43, meghna at index 44).
$ws.Range("A2").Value = "http://10.133.146.17:56080/SG/UOB_OCM"
$ws.Range("C2").Value = "P@ssw0rd@123"
$ws.Range("B2").Value = "meghna"

# Move the active selection from C2 to B3.
$ws.Activate()
$ws.Range("B3").Select()
